$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 1394.8928  # ALC!H33 (1442.2222 -> 1394.8928)
$ws.Cells.Item(33, 9).Value = 1508.9412  # ALC!I33 (1509 -> 1508.9412)
$ws.Cells.Item(33, 10).Value = 1218.6364  # ALC!J33 (1328.7 -> 1218.6364)
$ws.Cells.Item(33, 11).Value = 1508.9412  # ALC!K33 (1509 -> 1508.9412)
$ws.Cells.Item(33, 12).Value = 1218.6364  # ALC!L33 (1328.7 -> 1218.6364)
$ws.Cells.Item(33, 13).Value = -1279.9412  # ALC!M33 (-1280 -> -1279.9412)
$ws.Cells.Item(33, 14).Value = -1676.6364  # ALC!N33 (-1786.7 -> -1676.6364)

$ws.Cells.Item(53, 8).Value = 367.0909  # ALC!H53 (357.1 -> 367.0909)
$ws.Cells.Item(53, 10).Value = 486.125  # ALC!J53 (488.85715 -> 486.125)
$ws.Cells.Item(53, 12).Value = 486.125  # ALC!L53 (488.85715 -> 486.125)
$ws.Cells.Item(53, 14).Value = -1760.125  # ALC!N53 (-1762.85715 -> -1760.125)

$ws.Cells.Item(88, 8).Value = 2020224.8  # ALC!H88 (1616780.4 -> 2020224.8)
$ws.Cells.Item(88, 10).Value = 2692633  # ALC!J88 (2020225.5 -> 2692633)
$ws.Cells.Item(88, 12).Value = 2692633  # ALC!L88 (2020225.5 -> 2692633)
$ws.Cells.Item(88, 14).Value = -2693445  # ALC!N88 (-2021037.5 -> -2693445)

$ws.Cells.Item(91, 8).Value = 2020224.8  # ALC!H91 (1616780.4 -> 2020224.8)
$ws.Cells.Item(91, 10).Value = 2692633  # ALC!J91 (2020225.5 -> 2692633)
$ws.Cells.Item(91, 12).Value = 2692633  # ALC!L91 (2020225.5 -> 2692633)
$ws.Cells.Item(91, 14).Value = -2695441  # ALC!N91 (-2023033.5 -> -2695441)

$ws.Cells.Item(132, 8).Value = 2025.88  # ALC!H132 (2025.92 -> 2025.88)
$ws.Cells.Item(132, 9).Value = 1893.2727  # ALC!I132 (1964.762 -> 1893.2727)
$ws.Cells.Item(132, 10).Value = 2998.3333  # ALC!J132 (2347 -> 2998.3333)
$ws.Cells.Item(132, 11).Value = 5679.8181  # ALC!K132 (5894.286 -> 5679.8181)
$ws.Cells.Item(132, 12).Value = 8994.999899999999  # ALC!L132 (7041 -> 8994.999899999999)
$ws.Cells.Item(132, 13).Value = -3149.8181  # ALC!M132 (-3364.286 -> -3149.8181)
$ws.Cells.Item(132, 14).Value = -14054.9999  # ALC!N132 (-12101 -> -14054.9999)

$ws.Cells.Item(137, 8).Value = 3573.6572  # ALC!H137 (3404.4211 -> 3573.6572)
$ws.Cells.Item(137, 9).Value = 1741.7  # ALC!I137 (1617.7084 -> 1741.7)
$ws.Cells.Item(137, 10).Value = 6016.2666  # ALC!J137 (6467.357 -> 6016.2666)
$ws.Cells.Item(137, 11).Value = 5225.1  # ALC!K137 (4853.1252 -> 5225.1)
$ws.Cells.Item(137, 12).Value = 18048.7998  # ALC!L137 (19402.071 -> 18048.7998)
$ws.Cells.Item(137, 13).Value = -2675.1  # ALC!M137 (-2303.1252 -> -2675.1)
$ws.Cells.Item(137, 14).Value = -23148.7998  # ALC!N137 (-24502.071 -> -23148.7998)

$ws.Cells.Item(138, 8).Value = 2066.9656  # ALC!H138 (2064.0847 -> 2066.9656)
$ws.Cells.Item(138, 10).Value = 2566.1794  # ALC!J138 (2549.45 -> 2566.1794)
$ws.Cells.Item(138, 12).Value = 7698.5382  # ALC!L138 (7648.349999999999 -> 7698.5382)
$ws.Cells.Item(138, 14).Value = -17978.5382  # ALC!N138 (-17928.35 -> -17978.5382)

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(9, 8).Value = 800000  # ARM!H9 (799999.5 -> 800000)
$ws.Cells.Item(9, 9).Value = 800000  # ARM!I9 (799999.5 -> 800000)
$ws.Cells.Item(9, 11).Value = 800000  # ARM!K9 (799999.5 -> 800000)
$ws.Cells.Item(9, 13).Value = -799830  # ARM!M9 (-799829.5 -> -799830)

$ws.Cells.Item(20, 8).Value = 800000  # ARM!H20 (799999.5 -> 800000)
$ws.Cells.Item(20, 9).Value = 800000  # ARM!I20 (799999.5 -> 800000)
$ws.Cells.Item(20, 11).Value = 800000  # ARM!K20 (799999.5 -> 800000)
$ws.Cells.Item(20, 13).Value = -799730  # ARM!M20 (-799729.5 -> -799730)

$ws.Cells.Item(88, 8).Value = 1681.9  # ARM!H88 (1901.625 -> 1681.9)
$ws.Cells.Item(88, 9).Value = 1802  # ARM!I88 (2021.2 -> 1802)
$ws.Cells.Item(88, 10).Value = 1501.75  # ARM!J88 (1702.3334 -> 1501.75)
$ws.Cells.Item(88, 11).Value = 1802  # ARM!K88 (2021.2 -> 1802)
$ws.Cells.Item(88, 12).Value = 1501.75  # ARM!L88 (1702.3334 -> 1501.75)
$ws.Cells.Item(88, 13).Value = -1396  # ARM!M88 (-1615.2 -> -1396)
$ws.Cells.Item(88, 14).Value = -2313.75  # ARM!N88 (-2514.3334 -> -2313.75)

$ws.Cells.Item(91, 8).Value = 1681.9  # ARM!H91 (1901.625 -> 1681.9)
$ws.Cells.Item(91, 9).Value = 1802  # ARM!I91 (2021.2 -> 1802)
$ws.Cells.Item(91, 10).Value = 1501.75  # ARM!J91 (1702.3334 -> 1501.75)
$ws.Cells.Item(91, 11).Value = 1802  # ARM!K91 (2021.2 -> 1802)
$ws.Cells.Item(91, 12).Value = 1501.75  # ARM!L91 (1702.3334 -> 1501.75)
$ws.Cells.Item(91, 13).Value = -398  # ARM!M91 (-617.2 -> -398)
$ws.Cells.Item(91, 14).Value = -4309.75  # ARM!N91 (-4510.3334 -> -4309.75)

$ws.Cells.Item(94, 8).Value = 48499.5  # ARM!H94 (29249.75 -> 48499.5)
$ws.Cells.Item(94, 10).Value = 48499.5  # ARM!J94 (29249.75 -> 48499.5)
$ws.Cells.Item(94, 12).Value = 48499.5  # ARM!L94 (29249.75 -> 48499.5)
$ws.Cells.Item(94, 14).Value = -50301.5  # ARM!N94 (-31051.75 -> -50301.5)

$ws.Cells.Item(132, 8).Value = 11212.424  # ARM!H132 (11419.322 -> 11212.424)
$ws.Cells.Item(132, 9).Value = 8888.958000000001  # ARM!I132 (8969.272000000001 -> 8888.958000000001)
$ws.Cells.Item(132, 11).Value = 26666.874  # ARM!K132 (26907.816 -> 26666.874)
$ws.Cells.Item(132, 13).Value = -24136.874  # ARM!M132 (-24377.816 -> -24136.874)

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 4696.909  # BSM!H20 (3978.3845 -> 4696.909)
$ws.Cells.Item(20, 9).Value = 5138.857  # BSM!I20 (3737 -> 5138.857)
$ws.Cells.Item(20, 10).Value = 3923.5  # BSM!J20 (4783 -> 3923.5)
$ws.Cells.Item(20, 11).Value = 5138.857  # BSM!K20 (3737 -> 5138.857)
$ws.Cells.Item(20, 12).Value = 3923.5  # BSM!L20 (4783 -> 3923.5)
$ws.Cells.Item(20, 13).Value = -4891.857  # BSM!M20 (-3490 -> -4891.857)
$ws.Cells.Item(20, 14).Value = -4417.5  # BSM!N20 (-5277 -> -4417.5)

$ws.Cells.Item(40, 8).Value = 90461.2  # BSM!H40 (93079 -> 90461.2)
$ws.Cells.Item(40, 10).Value = 88719.25  # BSM!J40 (91629 -> 88719.25)
$ws.Cells.Item(40, 12).Value = 88719.25  # BSM!L40 (91629 -> 88719.25)
$ws.Cells.Item(40, 14).Value = -89249.25  # BSM!N40 (-92159 -> -89249.25)

$ws.Cells.Item(134, 8).Value = 74229.47  # BSM!H134 (53715.57 -> 74229.47)
$ws.Cells.Item(134, 9).Value = 3341.4  # BSM!I134 (2865.6667 -> 3341.4)
$ws.Cells.Item(134, 10).Value = 216005.6  # BSM!J134 (180840.33 -> 216005.6)
$ws.Cells.Item(134, 11).Value = 10024.2  # BSM!K134 (8597.000100000001 -> 10024.2)
$ws.Cells.Item(134, 12).Value = 648016.8  # BSM!L134 (542520.99 -> 648016.8)
$ws.Cells.Item(134, 13).Value = -7489.200000000001  # BSM!M134 (-6062.000100000001 -> -7489.200000000001)
$ws.Cells.Item(134, 14).Value = -653086.8  # BSM!N134 (-547590.99 -> -653086.8)

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(92, 8).Value = 48996.668  # CRP!H92 (70996.664 -> 48996.668)
$ws.Cells.Item(92, 10).Value = 48996.668  # CRP!J92 (70996.664 -> 48996.668)
$ws.Cells.Item(92, 12).Value = 48996.668  # CRP!L92 (70996.664 -> 48996.668)
$ws.Cells.Item(92, 14).Value = -53988.668  # CRP!N92 (-75988.664 -> -53988.668)

$ws.Cells.Item(122, 8).Value = 1411.2858  # CRP!H122 (1524.75 -> 1411.2858)
$ws.Cells.Item(122, 9).Value = 1322.25  # CRP!I122 (1499.6666 -> 1322.25)
$ws.Cells.Item(122, 10).Value = 1530  # CRP!J122 (1600 -> 1530)
$ws.Cells.Item(122, 11).Value = 3966.75  # CRP!K122 (4498.9998 -> 3966.75)
$ws.Cells.Item(122, 12).Value = 4590  # CRP!L122 (4800 -> 4590)
$ws.Cells.Item(122, 13).Value = -1516.75  # CRP!M122 (-2048.9998 -> -1516.75)
$ws.Cells.Item(122, 14).Value = -9490  # CRP!N122 (-9700 -> -9490)

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 453620.9  # CUL!H12 (476301.6 -> 453620.9)
$ws.Cells.Item(12, 9).Value = 2252.111  # CUL!I12 (2532.75 -> 2252.111)
$ws.Cells.Item(12, 11).Value = 6756.333  # CUL!K12 (7598.25 -> 6756.333)
$ws.Cells.Item(12, 13).Value = -6583.333  # CUL!M12 (-7425.25 -> -6583.333)

$ws.Cells.Item(45, 8).Value = 6367.6  # CUL!H45 (20291 -> 6367.6)
$ws.Cells.Item(45, 9).Value = 2920.5  # CUL!I45 (841 -> 2920.5)
$ws.Cells.Item(45, 10).Value = 8665.666999999999  # CUL!J45 (30016 -> 8665.666999999999)
$ws.Cells.Item(45, 11).Value = 8761.5  # CUL!K45 (2523 -> 8761.5)
$ws.Cells.Item(45, 12).Value = 25997.001  # CUL!L45 (90048 -> 25997.001)
$ws.Cells.Item(45, 13).Value = -8229.5  # CUL!M45 (-1991 -> -8229.5)
$ws.Cells.Item(45, 14).Value = -27061.001  # CUL!N45 (-91112 -> -27061.001)

$ws.Cells.Item(128, 8).Value = 446665  # CUL!H128 (419997.5 -> 446665)
$ws.Cells.Item(128, 9).Value = 446665  # CUL!I128 (419997.5 -> 446665)
$ws.Cells.Item(128, 11).Value = 1339995  # CUL!K128 (1259992.5 -> 1339995)
$ws.Cells.Item(128, 13).Value = -1335015  # CUL!M128 (-1255012.5 -> -1335015)

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 132.81818  # GSM!H2 (418.8889 -> 132.81818)
$ws.Cells.Item(2, 9).Value = 120.76923  # GSM!I2 (126.63636 -> 120.76923)
$ws.Cells.Item(2, 10).Value = 150.22223  # GSM!J2 (878.1429000000001 -> 150.22223)
$ws.Cells.Item(2, 11).Value = 120.76923  # GSM!K2 (126.63636 -> 120.76923)
$ws.Cells.Item(2, 12).Value = 150.22223  # GSM!L2 (878.1429000000001 -> 150.22223)
$ws.Cells.Item(2, 13).Value = -7.769229999999993  # GSM!M2 (-13.63636 -> -7.769229999999993)
$ws.Cells.Item(2, 14).Value = -376.22223  # GSM!N2 (-1104.1429 -> -376.22223)

$ws.Cells.Item(70, 8).Value = 3120.3333  # GSM!H70 (3358.5 -> 3120.3333)
$ws.Cells.Item(70, 9).Value = 3061.6365  # GSM!I70 (3608 -> 3061.6365)
$ws.Cells.Item(70, 10).Value = 3212.5715  # GSM!J70 (2942.6667 -> 3212.5715)
$ws.Cells.Item(70, 11).Value = 3061.6365  # GSM!K70 (3608 -> 3061.6365)
$ws.Cells.Item(70, 12).Value = 3212.5715  # GSM!L70 (2942.6667 -> 3212.5715)
$ws.Cells.Item(70, 13).Value = -2791.6365  # GSM!M70 (-3338 -> -2791.6365)
$ws.Cells.Item(70, 14).Value = -3752.5715  # GSM!N70 (-3482.6667 -> -3752.5715)

$ws.Cells.Item(73, 8).Value = 3120.3333  # GSM!H73 (3358.5 -> 3120.3333)
$ws.Cells.Item(73, 9).Value = 3061.6365  # GSM!I73 (3608 -> 3061.6365)
$ws.Cells.Item(73, 10).Value = 3212.5715  # GSM!J73 (2942.6667 -> 3212.5715)
$ws.Cells.Item(73, 11).Value = 3061.6365  # GSM!K73 (3608 -> 3061.6365)
$ws.Cells.Item(73, 12).Value = 3212.5715  # GSM!L73 (2942.6667 -> 3212.5715)
$ws.Cells.Item(73, 13).Value = -2125.6365  # GSM!M73 (-2672 -> -2125.6365)
$ws.Cells.Item(73, 14).Value = -5084.5715  # GSM!N73 (-4814.6667 -> -5084.5715)

$ws.Cells.Item(74, 8).Value = 0  # GSM!H74 (30000 -> 0)
$ws.Cells.Item(74, 10).Value = 0  # GSM!J74 (30000 -> 0)
$ws.Cells.Item(74, 12).Value = 0  # GSM!L74 (30000 -> 0)
$ws.Cells.Item(74, 14).Value = $null  # GSM!N74 clear (was -31872)

$ws.Cells.Item(77, 8).Value = 0  # GSM!H77 (30000 -> 0)
$ws.Cells.Item(77, 10).Value = 0  # GSM!J77 (30000 -> 0)
$ws.Cells.Item(77, 12).Value = 0  # GSM!L77 (90000 -> 0)
$ws.Cells.Item(77, 14).Value = $null  # GSM!N77 clear (was -99360)

$ws.Cells.Item(132, 8).Value = 45464370  # GSM!H132 (43492004 -> 45464370)
$ws.Cells.Item(132, 10).Value = 52505.5  # GSM!J132 (68336.664 -> 52505.5)
$ws.Cells.Item(132, 12).Value = 157516.5  # GSM!L132 (205009.992 -> 157516.5)
$ws.Cells.Item(132, 14).Value = -162576.5  # GSM!N132 (-210069.992 -> -162576.5)

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 3553.0356  # LTW!H40 (3442.8667 -> 3553.0356)
$ws.Cells.Item(40, 9).Value = 2899.3333  # LTW!I40 (2812.4783 -> 2899.3333)
$ws.Cells.Item(40, 11).Value = 2899.3333  # LTW!K40 (2812.4783 -> 2899.3333)
$ws.Cells.Item(40, 13).Value = -2763.3333  # LTW!M40 (-2676.4783 -> -2763.3333)

$ws.Cells.Item(46, 8).Value = 3999.913  # LTW!H46 (3980.3076 -> 3999.913)
$ws.Cells.Item(46, 10).Value = 4461.1113  # LTW!J46 (4303.3335 -> 4461.1113)
$ws.Cells.Item(46, 12).Value = 4461.1113  # LTW!L46 (4303.3335 -> 4461.1113)
$ws.Cells.Item(46, 14).Value = -4837.1113  # LTW!N46 (-4679.3335 -> -4837.1113)

$ws.Cells.Item(68, 8).Value = 2598.6667  # LTW!H68 (2749 -> 2598.6667)
$ws.Cells.Item(68, 9).Value = 2598.4  # LTW!I68 (2749 -> 2598.4)
$ws.Cells.Item(68, 10).Value = 2599.2  # LTW!J68 (0 -> 2599.2)
$ws.Cells.Item(68, 11).Value = 2598.4  # LTW!K68 (2749 -> 2598.4)
$ws.Cells.Item(68, 12).Value = 2599.2  # LTW!L68 (0 -> 2599.2)
$ws.Cells.Item(68, 13).Value = -1849.4  # LTW!M68 (-2000 -> -1849.4)
$ws.Cells.Item(68, 14).Value = -4097.2  # LTW!N68 (None -> -4097.2)

$ws.Cells.Item(71, 8).Value = 2598.6667  # LTW!H71 (2749 -> 2598.6667)
$ws.Cells.Item(71, 9).Value = 2598.4  # LTW!I71 (2749 -> 2598.4)
$ws.Cells.Item(71, 10).Value = 2599.2  # LTW!J71 (0 -> 2599.2)
$ws.Cells.Item(71, 11).Value = 12992  # LTW!K71 (13745 -> 12992)
$ws.Cells.Item(71, 12).Value = 12996  # LTW!L71 (0 -> 12996)
$ws.Cells.Item(71, 13).Value = -9248  # LTW!M71 (-10001 -> -9248)
$ws.Cells.Item(71, 14).Value = -20484  # LTW!N71 (None -> -20484)

$ws.Cells.Item(132, 8).Value = 79848.5  # LTW!H132 (105865.914 -> 79848.5)
$ws.Cells.Item(132, 9).Value = 10226.75  # LTW!I132 (14442 -> 10226.75)
$ws.Cells.Item(132, 11).Value = 30680.25  # LTW!K132 (43326 -> 30680.25)
$ws.Cells.Item(132, 13).Value = -28150.25  # LTW!M132 (-40796 -> -28150.25)

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(8, 8).Value = 3503.5  # WVR!H8 (4002.3333 -> 3503.5)
$ws.Cells.Item(8, 10).Value = 5004  # WVR!J8 (5002 -> 5004)
$ws.Cells.Item(8, 12).Value = 5004  # WVR!L8 (5002 -> 5004)
$ws.Cells.Item(8, 14).Value = -5284  # WVR!N8 (-5282 -> -5284)

$ws.Cells.Item(45, 8).Value = 8458.333000000001  # WVR!H45 (8491.666999999999 -> 8458.333000000001)
$ws.Cells.Item(45, 10).Value = 8458.333000000001  # WVR!J45 (8491.666999999999 -> 8458.333000000001)
$ws.Cells.Item(45, 12).Value = 8458.333000000001  # WVR!L45 (8491.666999999999 -> 8458.333000000001)
$ws.Cells.Item(45, 14).Value = -9440.333000000001  # WVR!N45 (-9473.666999999999 -> -9440.333000000001)

$ws.Cells.Item(136, 8).Value = 10721.28  # WVR!H136 (11185.708 -> 10721.28)
$ws.Cells.Item(136, 9).Value = 779.55554  # WVR!I136 (850.41174 -> 779.55554)
$ws.Cells.Item(136, 11).Value = 2338.66662  # WVR!K136 (2551.23522 -> 2338.66662)
$ws.Cells.Item(136, 13).Value = 211.33338  # WVR!M136 (-1.235220000000027 -> 211.33338)
